$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.999.19'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.67%  '
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("E12").Value = '  -1.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '3.043.33'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '62.919.60'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '2.581.74'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("E21").Value = '  -3.58%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  +3.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("E25").Value = '  +7.40%  '
$ws.Range("E26").Value = '  -2.42%  '
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  -2.67%  '
$ws.Range("E31").Value = '  -3.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '460.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '176.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0539'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.46%  '
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.06%  '
$ws.Range("E51").Value = '  +0.12%  '
